$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing monthly figures for 01-04-2021 (row 173) ---
$ws.Range("I173").Value = 7176
$ws.Range("K173").Value = 14243
$ws.Range("L173").Value = 57978

# --- Update existing monthly figures for 01-05-2021 (row 174) ---
$ws.Range("H174").Value = 23033
$ws.Range("I174").Value = 7911
$ws.Range("K174").Value = 14636
$ws.Range("L174").Value = 58952

# --- Update existing monthly figures for 01-06-2021 (row 175) ---
$ws.Range("G175").Value = 9736
$ws.Range("H175").Value = 15911
$ws.Range("I175").Value = 7227
$ws.Range("K175").Value = 14869
$ws.Range("L175").Value = 59846

# --- Append new monthly row for 01-07-2021 (row 176) ---
# A leading apostrophe forces Excel to store the date-like text as a
# text value (shared string) instead of auto-converting it to a date
# serial number; resetting the style afterwards keeps the cell's
# formatting identical to the other "Serie" cells in column A.
$ws.Range("A176").Value = "'01-07-2021"
$ws.Range("A176").Style = "Normal"

$ws.Range("B176").Value = 7325
$ws.Range("C176").Value = 64083
$ws.Range("D176").Value = 56758
$ws.Range("E176").Value = 120841
$ws.Range("F176").Value = 1495
$ws.Range("G176").Value = 5494
$ws.Range("H176").Value = 18430
$ws.Range("I176").Value = 7522
$ws.Range("J176").Value = 12157
$ws.Range("K176").Value = 14635
$ws.Range("L176").Value = 61109
